$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage so numeric-looking price strings (e.g. "309.86", "1.000")
# keep their exact original formatting instead of being parsed as numbers.
$ws.Range("B2:E51").NumberFormat = "@"

# Update cryptocurrency price/volume data per the latest GitHub Actions refresh
$ws.Range('D2').Value = '24.965.05'
$ws.Range('E2').Value = '  -3.12%  '
$ws.Range('D3').Value = '1.684.05'
$ws.Range('E3').Value = '  -2.67%  '
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  +0.75%  '
$ws.Range('D5').Value = '309.86'
$ws.Range('E5').Value = '  -0.93%  '
$ws.Range('D6').Value = '0.9962'
$ws.Range('E6').Value = '  +0.76%  '
$ws.Range('D7').Value = '0.3681'
$ws.Range('E7').Value = '  -2.10%  '
$ws.Range('D8').Value = '0.3376'
$ws.Range('E8').Value = '  -5.14%  '
$ws.Range('D9').Value = '47.66'
$ws.Range('E9').Value = '  -5.48%  '
$ws.Range('D10').Value = '1.181'
$ws.Range('E10').Value = '  -2.70%  '
$ws.Range('D11').Value = '0.07347'
$ws.Range('E11').Value = '  -2.23%  '
$ws.Range('D12').Value = '0.9971'
$ws.Range('E12').Value = '  +0.68%  '
$ws.Range('D13').Value = '6.206'
$ws.Range('E13').Value = '  -1.97%  '
$ws.Range('D14').Value = '20.61'
$ws.Range('E14').Value = '  -4.55%  '
$ws.Range('D15').Value = '6.859'
$ws.Range('E15').Value = '  -1.24%  '
$ws.Range('D16').Value = '1.681.59'
$ws.Range('E16').Value = '  -2.80%  '
$ws.Range('D17').Value = '0.00001105'
$ws.Range('E17').Value = '  -3.09%  '
$ws.Range('D18').Value = '0.06615'
$ws.Range('E18').Value = '  -2.16%  '
$ws.Range('D19').Value = '0.9963'
$ws.Range('E19').Value = '  +0.87%  '
$ws.Range('D20').Value = '82.52'
$ws.Range('E20').Value = '  -3.79%  '
$ws.Range('D21').Value = '16.93'
$ws.Range('E21').Value = '  -1.98%  '
$ws.Range('D22').Value = '6.215'
$ws.Range('E22').Value = '  -1.52%  '
$ws.Range('D23').Value = '12.66'
$ws.Range('E23').Value = '  +1.19%  '
$ws.Range('D24').Value = '24.863.16'
$ws.Range('E24').Value = '  -3.14%  '
$ws.Range('D25').Value = '2.428'
$ws.Range('E25').Value = '  +0.40%  '
$ws.Range('D26').Value = '2.718'
$ws.Range('E26').Value = '  -3.83%  '
$ws.Range('D27').Value = '19.91'
$ws.Range('E27').Value = '  -2.16%  '
$ws.Range('D28').Value = '150.91'
$ws.Range('E28').Value = '  -2.51%  '
$ws.Range('E29').Value = '  +10.97%  '
$ws.Range('D30').Value = '130.65'
$ws.Range('E30').Value = '  -2.18%  '
$ws.Range('D31').Value = '1.874.96'
$ws.Range('E31').Value = '  -2.35%  '
$ws.Range('D32').Value = '6.558'
$ws.Range('E32').Value = '  -1.28%  '
$ws.Range('D33').Value = '4.164'
$ws.Range('E33').Value = '  +1.77%  '
$ws.Range('D34').Value = '13.53'
$ws.Range('E34').Value = '  +1.28%  '
$ws.Range('B35').Value = 'Stellar'
$ws.Range('C35').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D35').Value = '0.08618'
$ws.Range('E35').Value = '  +0.93%  '
$ws.Range('B36').Value = 'WEMIXTOKEN'
$ws.Range('C36').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D36').Value = '1.733'
$ws.Range('E36').Value = '  -2.33%  '
$ws.Range('D37').Value = '5.463'
$ws.Range('E37').Value = '  -1.16%  '
$ws.Range('D38').Value = '0.06495'
$ws.Range('E38').Value = '  -2.30%  '
$ws.Range('D39').Value = '0.02359'
$ws.Range('E39').Value = '  -3.26%  '
$ws.Range('D40').Value = '8.782'
$ws.Range('E40').Value = '  -3.81%  '
$ws.Range('D41').Value = '0.2180'
$ws.Range('E41').Value = '  +0.33%  '
$ws.Range('D42').Value = '1.247'
$ws.Range('E42').Value = '  -2.09%  '
$ws.Range('D43').Value = '0.6299'
$ws.Range('E43').Value = '  -1.77%  '
$ws.Range('D44').Value = '0.9957'
$ws.Range('E44').Value = '  +0.76%  '
$ws.Range('D45').Value = '13.50'
$ws.Range('E45').Value = '  -0.86%  '
$ws.Range('D46').Value = '3.799'
$ws.Range('E46').Value = '  -1.70%  '
$ws.Range('D47').Value = '0.6000'
$ws.Range('E47').Value = '  -3.25%  '
$ws.Range('D48').Value = '2.055'
$ws.Range('E48').Value = '  -2.98%  '
$ws.Range('E49').Value = '  -3.79%  '
$ws.Range('D50').Value = '0.07180'
$ws.Range('E50').Value = '  -3.65%  '
$ws.Range('D51').Value = '77.68'
$ws.Range('E51').Value = '  -0.86%  '
